# Scheduled-runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures across the Sheets workbook (ALC, ARM, BSM, CRP, CUL, LTW, WVR tabs).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 113056.664
$ws.Range("I113").Value = 126938.75
$ws.Range("K113").Value = 126938.75
$ws.Range("M113").Value = -123684.75

$ws.Range("H123").Value = 17660.6
$ws.Range("J123").Value = 17660.6
$ws.Range("L123").Value = 17660.6
$ws.Range("N123").Value = -27460.6

$ws.Range("H137").Value = 2187.4
$ws.Range("I137").Value = 1575.7084
$ws.Range("J137").Value = 4634.1665
$ws.Range("K137").Value = 4727.1252
$ws.Range("L137").Value = 13902.4995
$ws.Range("M137").Value = -2177.1252
$ws.Range("N137").Value = -19002.4995

$ws.Range("H138").Value = 3479.6794
$ws.Range("I138").Value = 1106.3077
$ws.Range("J138").Value = 5853.0513
$ws.Range("K138").Value = 3318.9231
$ws.Range("L138").Value = 17559.1539
$ws.Range("M138").Value = 1821.0769
$ws.Range("N138").Value = -27839.1539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3425.64
$ws.Range("I45").Value = 2977.7856
$ws.Range("J45").Value = 3995.6365
$ws.Range("K45").Value = 2977.7856
$ws.Range("L45").Value = 3995.6365
$ws.Range("M45").Value = -2600.7856
$ws.Range("N45").Value = -4749.636500000001

$ws.Range("H61").Value = 2102.5
$ws.Range("I61").Value = 1014.7692
$ws.Range("K61").Value = 1014.7692
$ws.Range("M61").Value = -802.7692

$ws.Range("H88").Value = 2269.3157
$ws.Range("I88").Value = 2186.1428
$ws.Range("J88").Value = 2502.2
$ws.Range("K88").Value = 2186.1428
$ws.Range("L88").Value = 2502.2
$ws.Range("M88").Value = -1780.1428
$ws.Range("N88").Value = -3314.2

$ws.Range("H91").Value = 2269.3157
$ws.Range("I91").Value = 2186.1428
$ws.Range("J91").Value = 2502.2
$ws.Range("K91").Value = 2186.1428
$ws.Range("L91").Value = 2502.2
$ws.Range("M91").Value = -782.1428000000001
$ws.Range("N91").Value = -5310.2

$ws.Range("H132").Value = 2328.3948
$ws.Range("I132").Value = 2186.2666
$ws.Range("J132").Value = 2861.375
$ws.Range("K132").Value = 6558.7998
$ws.Range("L132").Value = 8584.125
$ws.Range("M132").Value = -4028.7998
$ws.Range("N132").Value = -13644.125

$ws.Range("H136").Value = 2102.5
$ws.Range("I136").Value = 1014.7692
$ws.Range("K136").Value = 3044.3076
$ws.Range("M136").Value = -494.3076000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 70533.56
$ws.Range("I86").Value = 158164.28
$ws.Range("J86").Value = 2376.3333
$ws.Range("K86").Value = 158164.28
$ws.Range("L86").Value = 2376.3333
$ws.Range("M86").Value = -157041.28
$ws.Range("N86").Value = -4622.3333

$ws.Range("H89").Value = 70533.56
$ws.Range("I89").Value = 158164.28
$ws.Range("J89").Value = 2376.3333
$ws.Range("K89").Value = 790821.4
$ws.Range("L89").Value = 11881.6665
$ws.Range("M89").Value = -785205.4
$ws.Range("N89").Value = -23113.6665

$ws.Range("H134").Value = 1997.1765
$ws.Range("I134").Value = 1892.0667
$ws.Range("J134").Value = 2785.5
$ws.Range("K134").Value = 5676.2001
$ws.Range("L134").Value = 8356.5
$ws.Range("M134").Value = -3141.2001
$ws.Range("N134").Value = -13426.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2083.86
$ws.Range("I31").Value = 904.4068
$ws.Range("J31").Value = 3781.122
$ws.Range("K31").Value = 904.4068
$ws.Range("L31").Value = 3781.122
$ws.Range("M31").Value = -609.4068
$ws.Range("N31").Value = -4371.121999999999

$ws.Range("H34").Value = 2083.86
$ws.Range("I34").Value = 904.4068
$ws.Range("J34").Value = 3781.122
$ws.Range("K34").Value = 904.4068
$ws.Range("L34").Value = 3781.122
$ws.Range("M34").Value = -702.4068
$ws.Range("N34").Value = -4185.121999999999

$ws.Range("H58").Value = 7757.641
$ws.Range("I58").Value = 1143.0322
$ws.Range("J58").Value = 33389.25
$ws.Range("K58").Value = 1143.0322
$ws.Range("L58").Value = 33389.25
$ws.Range("M58").Value = -940.0322000000001
$ws.Range("N58").Value = -33795.25

$ws.Range("H62").Value = 11113271
$ws.Range("J62").Value = 2700
$ws.Range("L62").Value = 2700
$ws.Range("N62").Value = -3948

$ws.Range("H65").Value = 11113271
$ws.Range("J65").Value = 2700
$ws.Range("L65").Value = 13500
$ws.Range("N65").Value = -19740

$ws.Range("H105").Value = 1266.6666
$ws.Range("I105").Value = 998
$ws.Range("K105").Value = 998
$ws.Range("M105").Value = 749

$ws.Range("H118").Value = 44990
$ws.Range("J118").Value = 44990
$ws.Range("L118").Value = 44990
$ws.Range("N118").Value = -48304

$ws.Range("H122").Value = 3470.75
$ws.Range("I122").Value = 3512.2666
$ws.Range("K122").Value = 10536.7998
$ws.Range("M122").Value = -8086.799800000001

$ws.Range("H134").Value = 1012.55554
$ws.Range("I134").Value = 1012.55554
$ws.Range("K134").Value = 3037.66662
$ws.Range("M134").Value = -502.66662

$ws.Range("H136").Value = 7757.641
$ws.Range("I136").Value = 1143.0322
$ws.Range("J136").Value = 33389.25
$ws.Range("K136").Value = 3429.0966
$ws.Range("L136").Value = 100167.75
$ws.Range("M136").Value = -879.0966000000003
$ws.Range("N136").Value = -105267.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6306.081
$ws.Range("I5").Value = 925.8929000000001
$ws.Range("J5").Value = 23044.445
$ws.Range("K5").Value = 2777.6787
$ws.Range("L5").Value = 69133.33499999999
$ws.Range("M5").Value = -2665.6787
$ws.Range("N5").Value = -69357.33499999999

$ws.Range("H23").Value = 701.86664
$ws.Range("I23").Value = 389.2
$ws.Range("J23").Value = 858.2
$ws.Range("K23").Value = 1167.6
$ws.Range("L23").Value = 2574.6
$ws.Range("M23").Value = -932.5999999999999
$ws.Range("N23").Value = -3044.6

$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16372

$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51864

$ws.Range("H69").Value = 1733.2222
$ws.Range("J69").Value = 1824.875
$ws.Range("L69").Value = 5474.625
$ws.Range("N69").Value = -7096.625

$ws.Range("H72").Value = 1733.2222
$ws.Range("J72").Value = 1824.875
$ws.Range("L72").Value = 16423.875
$ws.Range("N72").Value = -24535.875

$ws.Range("H80").Value = 12146.9
$ws.Range("I80").Value = 1695
$ws.Range("J80").Value = 13308.223
$ws.Range("K80").Value = 5085
$ws.Range("L80").Value = 39924.669
$ws.Range("M80").Value = -4149
$ws.Range("N80").Value = -41796.669

$ws.Range("H83").Value = 12146.9
$ws.Range("I83").Value = 1695
$ws.Range("J83").Value = 13308.223
$ws.Range("K83").Value = 15255
$ws.Range("L83").Value = 119774.007
$ws.Range("M83").Value = -10575
$ws.Range("N83").Value = -129134.007

$ws.Range("H113").Value = 617.4400000000001
$ws.Range("I113").Value = 566.25
$ws.Range("J113").Value = 708.44446
$ws.Range("K113").Value = 1698.75
$ws.Range("L113").Value = 2125.33338
$ws.Range("M113").Value = 471.25
$ws.Range("N113").Value = -6465.33338

$ws.Range("H120").Value = 4515
$ws.Range("I120").Value = 1030
$ws.Range("J120").Value = 8000
$ws.Range("K120").Value = 3090
$ws.Range("L120").Value = 24000
$ws.Range("M120").Value = 1748
$ws.Range("N120").Value = -33676

$ws.Range("H134").Value = 3060
$ws.Range("I134").Value = 3060
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9180
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4110
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 6306.081
$ws.Range("I135").Value = 925.8929000000001
$ws.Range("J135").Value = 23044.445
$ws.Range("K135").Value = 8333.036100000001
$ws.Range("L135").Value = 207400.005
$ws.Range("M135").Value = -5798.036100000001
$ws.Range("N135").Value = -212470.005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350

$ws.Range("H118").Value = 48396.668
$ws.Range("J118").Value = 48396.668
$ws.Range("L118").Value = 48396.668
$ws.Range("N118").Value = -51710.668

$ws.Range("H132").Value = 3829.8484
$ws.Range("I132").Value = 4390.7393
$ws.Range("J132").Value = 2539.8
$ws.Range("K132").Value = 13172.2179
$ws.Range("L132").Value = 7619.400000000001
$ws.Range("M132").Value = -10642.2179
$ws.Range("N132").Value = -12679.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 49985
$ws.Range("J116").Value = 49985
$ws.Range("L116").Value = 49985
$ws.Range("N116").Value = -59163

$ws.Range("H132").Value = 2243.7058
$ws.Range("I132").Value = 1911.0769
$ws.Range("J132").Value = 3324.75
$ws.Range("K132").Value = 5733.2307
$ws.Range("L132").Value = 9974.25
$ws.Range("M132").Value = -3203.2307
$ws.Range("N132").Value = -15034.25

$ws.Range("H136").Value = 1396.1842
$ws.Range("J136").Value = 2108.524
$ws.Range("L136").Value = 6325.572
$ws.Range("N136").Value = -11425.572
